$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "Highs-BigM (100,100)"
$ws.Range("B5").Value = "FEASIBLE_POINT"
$ws.Range("C5").Value = "OPTIMAL"
$ws.Range("D5").Value = -13.0
$ws.Range("E5").Value = 0.001993264
$ws.Range("F5").Value = 0.01565619
$ws.Range("G5").Value = 0.0035683174604316543
$ws.Range("H5").Value = 2597
$ws.Range("I5").Value = 0.116912841796875
$ws.Range("J5").Value = 5.0
$ws.Range("K5").Value = 4.0
$ws.Range("L5").Value = 1.9999999999999996

# Row 6
$ws.Range("A6").Value = "SOS1"
$ws.Range("B6").Value = "FEASIBLE_POINT"
$ws.Range("C6").Value = "OPTIMAL"
$ws.Range("D6").Value = -13.0
$ws.Range("E6").Value = 0.004863592
$ws.Range("F6").Value = 0.02544338
$ws.Range("G6").Value = 0.009506216718095237
$ws.Range("H6").Value = 3079
$ws.Range("I6").Value = 0.1317291259765625
$ws.Range("J6").Value = 5.0
$ws.Range("K6").Value = 4.0
$ws.Range("L6").Value = 2.0

# Row 7
$ws.Range("A7").Value = "Product_Mode"
$ws.Range("B7").Value = "FEASIBLE_POINT"
$ws.Range("C7").Value = "LOCALLY_SOLVED"
$ws.Range("D7").Value = -13.000000246414976
$ws.Range("E7").Value = 0.01481231
$ws.Range("F7").Value = 0.032895627
$ws.Range("G7").Value = 0.020900470322175735
$ws.Range("H7").Value = 10616
$ws.Range("I7").Value = 0.29681396484375
$ws.Range("J7").Value = 5.000000032830815
$ws.Range("K7").Value = 4.000000038995717
$ws.Range("L7").Value = 1.9999999517014941
